$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the old "pdfs:" / "2013, 2015-2021" header row and the lone 2014 value
$ws.Range("A1:B2").Clear() | Out-Null

# Years for rows 1-8 (row 9 holds 2014 separately, added below)
$years = @(2013, 2015, 2016, 2017, 2018, 2019, 2020, 2021)

for ($i = 0; $i -lt $years.Length; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 1).Value = $years[$i]
    $ws.Cells.Item($row, 2).Formula = '="bills_laws/lawsstatutes/"&A' + $row + '&"orlaw$$$$.pdf"'
}

# Row 9: 2014 uses a differently formatted literal PDF link (not a formula)
$ws.Cells.Item(9, 1).Value = 2014
$ws.Cells.Item(9, 2).Value = 'bills_laws/lawsstatutes/2014R1orLaw$$$$ss.pdf'

$ws.Range("B10").Select() | Out-Null
